# Corrected data cleaning for pre/post/total fixation data
#
# 1) Header row (row 1) loses the bold/centered/bordered "header" style —
#    it goes back to the default (unstyled) look, and the first header
#    cell (A1, "Unnamed: 0") is cleared out entirely.
# 2) A bunch of the computed metric cells in rows 3-8 are corrected, and
#    the "assign2" (G) / "var3" (AK) columns are blanked out for the
#    metric rows (3-8) since they no longer contain valid data.
# 3) The two trailing fully-blank rows (10 and 11) are removed, shrinking
#    the used range from A1:AP11 down to A1:AP9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Strip the bold/border/center styling back off the header row ---
$headerRange = $ws.Range("A1:AP1")
$headerRange.ClearFormats()

# The old "Unnamed: 0" header label is removed.
$ws.Range("A1").ClearContents()

# --- 2) Corrected metric values ---

# Row 3 - Revisit count
$ws.Range("B3").Value = 2
$ws.Range("G3").ClearContents()
$ws.Range("H3").Value = 38
$ws.Range("I3").Value = 27
$ws.Range("J3").Value = 7
$ws.Range("L3").Value = 5
$ws.Range("S3").Value = 20
$ws.Range("T3").Value = 10
$ws.Range("W3").Value = 0
$ws.Range("Y3").Value = 1
$ws.Range("Z3").Value = 16
$ws.Range("AB3").Value = 24
$ws.Range("AC3").Value = 5
$ws.Range("AK3").ClearContents()

# Row 4 - Fixation count
$ws.Range("B4").Value = 6
$ws.Range("G4").ClearContents()
$ws.Range("H4").Value = 115
$ws.Range("I4").Value = 59
$ws.Range("J4").Value = 16
$ws.Range("L4").Value = 9
$ws.Range("S4").Value = 39
$ws.Range("T4").Value = 24
$ws.Range("W4").Value = 1
$ws.Range("Y4").Value = 3
$ws.Range("Z4").Value = 25
$ws.Range("AB4").Value = 48
$ws.Range("AC4").Value = 6
$ws.Range("AK4").ClearContents()

# Row 5 - Dwell time (ms)
$ws.Range("B5").Value = 2127.19
$ws.Range("G5").ClearContents()
$ws.Range("H5").Value = 33834.31
$ws.Range("I5").Value = 17935.75
$ws.Range("J5").Value = 5823.06
$ws.Range("L5").Value = 4004.6
$ws.Range("S5").Value = 12387.93
$ws.Range("T5").Value = 8225.59
$ws.Range("W5").Value = 834.42
$ws.Range("Y5").Value = 1134.5
$ws.Range("Z5").Value = 8425.28
$ws.Range("AB5").Value = 13170.47
$ws.Range("AC5").Value = 1918.67
$ws.Range("AK5").ClearContents()

# Row 6 - Dwell time (%)
$ws.Range("B6").Value = 1.07
$ws.Range("C6").Value = 0.18
$ws.Range("F6").Value = 0.07
$ws.Range("G6").ClearContents()
$ws.Range("H6").Value = 17.04
$ws.Range("I6").Value = 9.03
$ws.Range("J6").Value = 2.93
$ws.Range("K6").Value = 2.34
$ws.Range("L6").Value = 2.02
$ws.Range("M6").Value = 0.2
$ws.Range("Q6").Value = 0.07
$ws.Range("S6").Value = 6.24
$ws.Range("T6").Value = 4.14
$ws.Range("U6").Value = 1.51
$ws.Range("V6").Value = 1.06
$ws.Range("W6").Value = 0.42
$ws.Range("Y6").Value = 0.57
$ws.Range("Z6").Value = 4.24
$ws.Range("AA6").Value = 0.56
$ws.Range("AB6").Value = 6.63
$ws.Range("AC6").Value = 0.97
$ws.Range("AJ6").Value = 0.3
$ws.Range("AK6").ClearContents()
$ws.Range("AM6").Value = 0.07
$ws.Range("AO6").Value = 0.14

# Row 7 - Fixation duration (ms)
$ws.Range("B7").Value = 354.53
$ws.Range("G7").ClearContents()
$ws.Range("H7").Value = 294.21
$ws.Range("I7").Value = 304
$ws.Range("J7").Value = 363.94
$ws.Range("L7").Value = 444.96
$ws.Range("S7").Value = 317.64
$ws.Range("T7").Value = 342.73
$ws.Range("W7").Value = 834.42
$ws.Range("Y7").Value = 378.17
$ws.Range("Z7").Value = 337.01
$ws.Range("AB7").Value = 274.38
$ws.Range("AC7").Value = 319.78
$ws.Range("AK7").ClearContents()

# Row 8 - First fixation duration (ms)
$ws.Range("G8").ClearContents()
$ws.Range("AK8").ClearContents()

# --- 3) Drop the two trailing blank rows (10 & 11) ---
$ws.Range("A10:AP11").EntireRow.Delete()
